$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BF2:BF31 hold a "Date" text column that was mis-derived (off by a day
# because of how the NBA stats were dated). Correct the text value for
# every data row, forcing it to stay plain text (not get auto-converted
# to a date serial by the COM layer).
$rng = $ws.Range("BF2:BF31")
$rng.NumberFormat = "@"
for ($row = 2; $row -le 31; $row++) {
    $ws.Range("BF$row").Value = "2014-05-26"
}
$rng.Style = "Normal"
